# StateFunction.xlsx - update SLG building state-function config
# Most rows had all "enabled" (1) flags for columns C..N; this edit disables
# (sets to 0) the function flags that are no longer applicable per row,
# leaving the first couple of columns (B, and a few row-specific ones) and
# the final column O (EFT_FINISH) untouched at 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2:N2").Value = 0

$ws.Range("C3:D3").Value = 0
$ws.Range("F3:N3").Value = 0

$ws.Range("C4").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("G4:N4").Value = 0

$ws.Range("C5:D5").Value = 0
$ws.Range("F5:N5").Value = 0

$ws.Range("C6:N6").Value = 0
$ws.Range("C7:N7").Value = 0
$ws.Range("C8:N8").Value = 0
$ws.Range("C9:N9").Value = 0
$ws.Range("C10:N10").Value = 0
$ws.Range("C11:N11").Value = 0
$ws.Range("C12:N12").Value = 0
$ws.Range("C13:N13").Value = 0

# Move active selection to F11 (matches the selection saved in the file)
$ws.Range("F11").Select()
